# Update the "startsAt" / "presentBetween" placeholder templates in the
# schedule header rows: the bracketed single-argument notation
# (e.g. "${startsAt: [5]}") is replaced with the new bracket-less,
# comma-separated call notation (e.g. "${startsAt: 5}" /
# "${presentBetween: 17, 18}").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4:AF4").Value = '${startsAt: 5}'
$ws.Range("B5:AF5").Value = '${startsAt: 7}'
$ws.Range("B6:AF6").Value = '${startsAt: 9}'
$ws.Range("B7:AF7").Value = '${startsAt: 11}'
$ws.Range("B8:AF8").Value = '${presentBetween: 17, 18}'
$ws.Range("B9:AF9").Value = '${startsAt: 19}'

# Move the active selection from Y13 to B4, matching the saved view state.
$ws.Range("B4").Select()
